# Outstandings.xlsx edit script
# Summary of changes:
#  1. Insert a new data row (row 22) on "Purchase 22-23" into the running
#     outstanding-balance block (rows 13-21), pushing all subsequent rows
#     down by one.
#  2. Fill the new row 22 with the new ledger entry (date 44998,
#     invoice "483/22-23", client "Namrata Rubber Product", amount 26621,
#     running-total formula).
#  3. Un-bold the old block-closing row (now row 21) and bold the new
#     block-closing row (row 22), matching the running-total block styling.
#     (The row insert already shifts every row below 22 - including the
#     "27000 * 18%" GST block and its SUM/total formulas - down by one and
#     keeps their relative formulas correct, so no further repair is
#     needed there.)
#  4. Switch the active sheet/tab selection from "Sale 22-23" back to
#     "Purchase 22-23", with a new selected cell (G23 on sheet1, and the
#     previous A13 selection remains untouched on sheet2, only losing the
#     tabSelected flag).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Purchase 22-23"
$ws2 = $wb.Worksheets.Item(2)   # "Sale 22-23"

# --- 1. Insert the new row -------------------------------------------------
$ws1.Rows(22).Insert()

# --- 2. Populate the new row 22 --------------------------------------------
$ws1.Range("B22").Value = 44998
$ws1.Range("C22").Value = "483/22-23"
$ws1.Range("D22").Value = "Namrata Rubber Product"
$ws1.Range("E22").Value = 26621
$ws1.Range("F22").Formula = "=F21+E22"

# --- 3. Fix up bold styling on the block-closing cell ----------------------
# Row 21 is no longer the last row of the running-total block -> un-bold.
$ws1.Range("F21").Font.Bold = $false
# Row 22 is now the last row of the running-total block -> bold.
$ws1.Range("F22").Font.Bold = $true

# --- 4. Switch active sheet / selection ------------------------------------
$ws1.Activate()
$ws1.Range("G23").Select()
